$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 (Uke 10): 09.03 Kontakttime / 11.03 Oversiktsforelesning: Hypotesetesting
$ws.Range("C9").Value = "09.03: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."
$ws.Range("D9").Value = "11.03: **Oversiktsforelesning: Hypotesetesting** på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."

# Row 10 (Uke 11): 16.03 Kontakttime / 18.03 Oppgaveseminar
$ws.Range("C10").Value = "16.03: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."
$ws.Range("D10").Value = "18.03: Oppgaveseminar på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09). Se \@ref(seminar) for oppgaver."

# Row 11 (Uke 12): 23.03 Kontakttime / 25.03 Oversiktsforelesning: Regresjon
$ws.Range("C11").Value = "23.03: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."
$ws.Range("D11").Value = "25.03: **Oversiktsforelesning: Regresjon** på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."

# Row 13 (Uke 14): 06.04 Kontakttime / 08.04 Oppgaveseminar
$ws.Range("C13").Value = "06.04: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."
$ws.Range("D13").Value = "08.04: Oppgaveseminar på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09). Se \@ref(seminar) for oppgaver."

# Row 14 (Uke 15): 13.04 Kontakttime
$ws.Range("C14").Value = "13.04: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."

# Update the last active selection to reflect the user's latest interaction
$ws.Range("E22").Select()
